# Update a handful of already-reported days with corrected source figures
# and enter the newly published day (2021-04-16, row 416).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections for previously entered days ---
$ws.Range("C410").Value = 82
$ws.Range("C413").Value = 107
$ws.Range("C414").Value = 145
$ws.Range("L414").Value = 1
$ws.Range("C415").Value = 92
$ws.Range("L415").Value = 0
$ws.Range("M415").Value = 0

# --- New day entered: 2021-04-16 ---
$ws.Range("C416").Value = 10
$ws.Range("E416").Value = 11
$ws.Range("F416").Value = 9
$ws.Range("G416").Value = 26
$ws.Range("L416").Value = 0
$ws.Range("M416").Value = 0

# --- View state: scrolled down to keep the new row in view ---
$ws.Range("O14").Select()
$ws.Application.ActiveWindow.ScrollRow = 414
